{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that contains the \"git add (for all relevant files)\" bullet.\nlet target = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.indexOf(\"git add (for all relevant files)\") !== -1) {\n    target = paras.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\n// Append \" (or \" as a plain run.\ntarget.insertText(\" (or \", Word.InsertLocation.end);\n\n// Append \"git add \u2013all\" as a Courier-formatted run. insertOoxml lets us set\n// the w:cs (complex-script) font alongside ascii/hAnsi, which the Range.font\n// object does not expose.\nconst runOoxml =\n  '<?xml version=\"1.0\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Courier\" w:hAnsi=\"Courier\" w:cs=\"Courier New\"/></w:rPr>' +\n  '<w:t>git add \u2013all</w:t></w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\ntarget.insertOoxml(runOoxml, Word.InsertLocation.end);\n\n// Append the closing paren as a plain run.\ntarget.insertText(\")\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the \"git add (for all relevant files)\" bullet.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*git add (for all relevant files)*\") {\n        $target = $p\n        break\n    }\n}\n\n$r = $target.Range\n\n# Append the three new runs at the end of the paragraph (before the paragraph mark).\n$r.InsertAfter(\" (or \")\n$r.InsertAfter(\"git add \u2013all\")\n$r.InsertAfter(\")\")\n\n# Apply the Courier / Courier New formatting to just the newly-inserted\n# \"git add \u2013all\" run, using Find/Replace with formatting so the change\n# stays scoped to that single run instead of leaking to the whole paragraph.\n$searchText = \"git add \u2013all\"\n$fr = $target.Range\n$fr.Find.ClearFormatting()\n$fr.Find.Replacement.ClearFormatting()\n$fr.Find.Replacement.Font.Name = \"Courier\"\n$fr.Find.Replacement.Font.NameBi = \"Courier New\"\n$fr.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $searchText, 1, $true)\n"}
